# sentiment_history.xlsx: switch data source naming from CSV-era sheet
# names (with spaces) to SQL-storage-era names (underscored), and turn
# the "Graphs" sheet into a "Data Analysis" sheet that pulls a live date
# value out of the DOW news data via formula.

$wb = $excel.ActiveWorkbook

# --- Rename worksheets (CSV-style display names -> sql_style names) ---
$wb.Worksheets.Item("DOW News").Name       = "DOW_news"
$wb.Worksheets.Item("DOW Twitter").Name    = "DOW_twitter"
$wb.Worksheets.Item("NASDAQ Twitter").Name = "NASDAQ_twitter"
$wb.Worksheets.Item("NASDAQ News").Name    = "NASDAQ_news"
$wb.Worksheets.Item("Graphs").Name         = "Data Analysis"

# --- Clear the stray/stale cell selection left on DOW_twitter ---
$twitter = $wb.Worksheets.Item("DOW_twitter")
$twitter.Range("A1").Select() | Out-Null

# --- Populate the "Data Analysis" sheet with a formula that references
#     the renamed DOW_news sheet, formatted with a 3-decimal number
#     format (new custom numFmt / cellXf) ---
$analysis = $wb.Worksheets.Item("Data Analysis")
$analysis.Range("A1").Formula = "=DOW_news!A$2"
$analysis.Range("A1").NumberFormat = "0.000"
$analysis.Columns.Item(1).EntireColumn.AutoFit() | Out-Null

# --- Put the focus back on the Data Analysis sheet (it was, and stays,
#     the active tab) ---
$analysis.Activate() | Out-Null
$analysis.Range("A1").Select() | Out-Null
